$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$desc = "L110 L111 L120 L121 L130 L132 L210 L211 L220 L222 L300 L301 L303 L310 L312 L313 L335 L350 L351 L353 L355 L358 L360 L361 L362 L363 L364 L365 L366 L375 L380 L381 L382 L383 L385 L386 L395 L396 L400 L405 L455 L456 L475 L485 L486 ME10 ME101 ME303 ME401 ET2500 ET2550 ET2600 ET2650 TX235W TX430W SX230 SX235W SX430W SX435W SX440W SX445W XP30 XP33 XP100 XP102 XP103 XP200 XP201 XP202 XP203 XP204 XP205 XP206 XP207 XP208 XP210 XP211 XP212 XP213 XP214 XP215 XP216 XP217 XP225 XP231 XP235 XP240 XP241 XP242 XP243 XP245 XP247 XP300 XP301 XP302 XP303 XP305 XP306 XP310 XP312 XP313 XP315 XP322 XP323 XP325 XP330 XP332 XP333 XP335 XP340 XP343 XP344 XP345 XP400 XP401 XP402 XP403 XP405 XP406 XP410 XP411 XP412 XP413 XP415 XP420 XP422 XP423 XP425 XP430 XP431 XP432 XP434 XP435 XP440 XP441 XP442 XP445 XP446 XP452 XP2100 XP2105 XP2150 XP2155"

$ws.Range("A3").Value = "INO44L"
$ws.Range("B3").Value = "Almohadilla Epson"
$ws.Range("C3").Value = $desc
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 100000
$ws.Range("F3").Value = 10
$ws.Range("G3").Value = 32
$ws.Range("H3").Formula = "=(E3-D3)*G3"
$ws.Range("I3").Formula = "=D3*F3"
$ws.Range("J3").Value = 0
